$p = $ppt.ActivePresentation
$d2 = $p.Designs.Add()
Write-Host "New design: $($d2.Name)"
Write-Host "New master index: $($d2.Index)"
